$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2
$ws.Range("G3").Value = 2
$ws.Range("G6").Value = 1
$ws.Range("G8").Value = 1
